$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '42.891.23'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '2.311.98'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('D5').Value = "'301.77"
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = "'96.73"
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').Value = "'0.505"
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = "'0.498"
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').Value = "'34.98"
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  +6.27%  '
$ws.Range('D12').Value = "'0.0791"
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '2.672.38'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '2.314.62'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = "'0.785"
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '42.848.00'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = "'12.45"
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = '0.0₃0891'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').Value = "'6.03"
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').Value = "'67.66"
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = "'235.50"
$ws.Range('D24').Value = "'2.24"
$ws.Range('E24').Value = '  +3.99%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = "'2.42"
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').Value = "'24.68"
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').Value = "'2.31"
$ws.Range('E28').Value = '  +12.51%  '
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').Value = "'9.07"
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = "'32.54"
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').Value = "'17.75"
$ws.Range('E34').Value = '  +3.25%  '
$ws.Range('D35').Value = "'4.49"
$ws.Range('E35').Value = '  -6.81%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.0698"
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = "'2.35"
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').Value = "'2.76"
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').Value = "'0.109"
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('D42').Value = '1.978.17'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').Value = "'10.59"
$ws.Range('E43').Value = '  +5.83%  '
$ws.Range('D44').Value = "'18.76"
$ws.Range('E44').Value = '  +5.06%  '
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('E46').Value = '  -2.53%  '
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').Value = "'53.28"
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('D51').Value = "'72.05"
$ws.Range('E51').Value = '  +0.16%  '
